$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.026409909042907
$ws.Range("D2").Value = 1.035091042317973
$ws.Range("E2").Value = 1.030039241472381
$ws.Range("F2").Value = 1.04386335903152
$ws.Range("I2").Value = 1.033178144066372
$ws.Range("J2").Value = 1.031573507755817
$ws.Range("K2").Value = 1.037888642323885
$ws.Range("L2").Value = 1.032851408866051
$ws.Range("M2").Value = 1.046636026768367
$ws.Range("N2").Value = 1.014534486990944

# Row 3
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.027306469732585
$ws.Range("D3").Value = 1.035754859373252
$ws.Range("E3").Value = 1.030883158265832
$ws.Range("F3").Value = 1.04466648462613
$ws.Range("I3").Value = 1.033309341724197
$ws.Range("J3").Value = 1.032110088853492
$ws.Range("K3").Value = 1.038362169641206
$ws.Range("L3").Value = 1.033503501140152
$ws.Range("M3").Value = 1.047250294595473
$ws.Range("N3").Value = 1.014714506794011

# Row 4
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.027887271470443
$ws.Range("D4").Value = 1.036184790439941
$ws.Range("E4").Value = 1.031430246615142
$ws.Range("F4").Value = 1.045186892443338
$ws.Range("I4").Value = 1.03339306787626
$ws.Range("J4").Value = 1.032457338243665
$ws.Range("K4").Value = 1.038668287911581
$ws.Range("L4").Value = 1.033925822948246
$ws.Range("M4").Value = 1.047647844302253
$ws.Range("N4").Value = 1.014830943175951

# Row 5
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.028131598846378
$ws.Range("D5").Value = 1.036365626816863
$ws.Range("E5").Value = 1.031660484544147
$ws.Range("F5").Value = 1.045405845264041
$ws.Range("I5").Value = 1.033427986209251
$ws.Range("J5").Value = 1.032603331604044
$ws.Range("K5").Value = 1.038796910412612
$ws.Range("L5").Value = 1.034103455317522
$ws.Range("M5").Value = 1.047814991203749
$ws.Range("N5").Value = 1.014879880973006

# Row 6
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.028172631732706
$ws.Range("D6").Value = 1.036395995465377
$ws.Range("E6").Value = 1.031699156647365
$ws.Range("F6").Value = 1.045442618519534
$ws.Range("I6").Value = 1.033433832704857
$ws.Range("J6").Value = 1.032627845066999
$ws.Range("K6").Value = 1.038818502577496
$ws.Range("L6").Value = 1.034133285710822
$ws.Range("M6").Value = 1.04784305685719
$ws.Range("N6").Value = 1.01488809712417

# Row 7
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.027890535564382
$ws.Range("D7").Value = 1.036187206419767
$ws.Range("E7").Value = 1.031433322116862
$ws.Range("F7").Value = 1.045189817422704
$ws.Range("I7").Value = 1.033393535558293
$ws.Range("J7").Value = 1.032459288977397
$ws.Range("K7").Value = 1.038670006846867
$ws.Range("L7").Value = 1.033928196135497
$ws.Range("M7").Value = 1.047650077662414
$ws.Range("N7").Value = 1.014831597133139

# Row 8
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.026712766995123
$ws.Range("D8").Value = 1.035315299169447
$ws.Range("E8").Value = 1.030324235369896
$ws.Range("F8").Value = 1.044134626448891
$ws.Range("I8").Value = 1.0332227242847
$ws.Range("J8").Value = 1.03175483770751
$ws.Range("K8").Value = 1.038048731472149
$ws.Range("L8").Value = 1.033071708132492
$ws.Range("M8").Value = 1.046843604531437
$ws.Range("N8").Value = 1.014595335288627

# Row 9
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.024642557111042
$ws.Range("D9").Value = 1.033782003016323
$ws.Range("E9").Value = 1.028377751215857
$ws.Range("F9").Value = 1.042280931267871
$ws.Range("I9").Value = 1.032912819532516
$ws.Range("J9").Value = 1.030513905751888
$ws.Range("K9").Value = 1.03695183320209
$ws.Range("L9").Value = 1.031565395894436
$ws.Range("M9").Value = 1.045423155058899
$ws.Range("N9").Value = 1.014178660600265

# Row 10
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.023265973036195
$ws.Range("D10").Value = 1.032762001472766
$ws.Range("E10").Value = 1.027085484444111
$ws.Range("F10").Value = 1.041049066632327
$ws.Range("I10").Value = 1.032700261464613
$ws.Range("J10").Value = 1.029686957962527
$ws.Range("K10").Value = 1.036219209066395
$ws.Range("L10").Value = 1.030563238844331
$ws.Range("M10").Value = 1.044476721953571
$ws.Range("N10").Value = 1.013900668398148

# Row 11
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.022670757142992
$ws.Range("D11").Value = 1.032320871667393
$ws.Range("E11").Value = 1.026527217775891
$ws.Range("F11").Value = 1.040516611285781
$ws.Range("I11").Value = 1.032606817032283
$ws.Range("J11").Value = 1.029328975911317
$ws.Range("K11").Value = 1.035901668165862
$ws.Range("L11").Value = 1.030129797451905
$ws.Range("M11").Value = 1.044067051375537
$ws.Range("N11").Value = 1.013780250789118

# Row 12
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.022449797106156
$ws.Range("D12").Value = 1.032157098900333
$ws.Range("E12").Value = 1.026320048670734
$ws.Range("F12").Value = 1.040318978268269
$ws.Range("I12").Value = 1.032571897089926
$ws.Range("J12").Value = 1.029196020409796
$ws.Range("K12").Value = 1.035783674043233
$ws.Range("L12").Value = 1.029968874398547
$ws.Range("M12").Value = 1.043914904214712
$ws.Range("N12").Value = 1.013735515997402

# Row 13
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.022497187895885
$ws.Range("D13").Value = 1.032192224954696
$ws.Range("E13").Value = 1.026364478251038
$ws.Range("F13").Value = 1.040361364647941
$ws.Range("I13").Value = 1.032579397052527
$ws.Range("J13").Value = 1.029224539119762
$ws.Range("K13").Value = 1.035808986212603
$ws.Range("L13").Value = 1.030003389461359
$ws.Range("M13").Value = 1.043947539244119
$ws.Range("N13").Value = 1.013745112041257

# Row 14
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.022652489864933
$ws.Range("D14").Value = 1.032307332466812
$ws.Range("E14").Value = 1.026510089105145
$ws.Range("F14").Value = 1.040500271918272
$ws.Range("I14").Value = 1.032603934829877
$ws.Range("J14").Value = 1.029317985455847
$ws.Range("K14").Value = 1.03589191565161
$ws.Range("L14").Value = 1.030116493936049
$ws.Range("M14").Value = 1.044054474371246
$ws.Range("N14").Value = 1.013776553123791

# Row 15
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.02274819375314
$ws.Range("D15").Value = 1.032378264985286
$ws.Range("E15").Value = 1.026599830779196
$ws.Range("F15").Value = 1.040585876489135
$ws.Range("I15").Value = 1.032619025487147
$ws.Range("J15").Value = 1.029375562846791
$ws.Range("K15").Value = 1.035943005251615
$ws.Range("L15").Value = 1.030186191497928
$ws.Range("M15").Value = 1.044120363690055
$ws.Range("N15").Value = 1.013795924187179

# Row 16
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.023305493662928
$ws.Range("D16").Value = 1.032791289278179
$ws.Range("E16").Value = 1.027122562157326
$ws.Range("F16").Value = 1.041084424117294
$ws.Range("I16").Value = 1.032706433507182
$ws.Range("J16").Value = 1.029710718087056
$ws.Range("K16").Value = 1.036240276789637
$ws.Range("L16").Value = 1.030592015557551
$ws.Range("M16").Value = 1.044503913545201
$ws.Range("N16").Value = 1.013908659205692

# Row 17
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.02365530277541
$ws.Range("D17").Value = 1.033050513646742
$ws.Range("E17").Value = 1.027450805196967
$ws.Range("F17").Value = 1.041397405569129
$ws.Range("I17").Value = 1.032760886427045
$ws.Range("J17").Value = 1.029920977524286
$ws.Range("K17").Value = 1.03642666530958
$ws.Range("L17").Value = 1.030846712895934
$ws.Range("M17").Value = 1.044744543111902
$ws.Range("N17").Value = 1.013979363136245

# Row 18
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.023859422830791
$ws.Range("D18").Value = 1.033201766590211
$ws.Range("E18").Value = 1.027642388487612
$ws.Range("F18").Value = 1.041580053908096
$ws.Range("I18").Value = 1.032792512252632
$ws.Range("J18").Value = 1.030043627074791
$ws.Range("K18").Value = 1.036535352603146
$ws.Range("L18").Value = 1.030995321619779
$ws.Range("M18").Value = 1.044884911789165
$ws.Range("N18").Value = 1.014020599155913

# Row 19
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.023929036416982
$ws.Range("D19").Value = 1.033253348661504
$ws.Range("E19").Value = 1.027707734535463
$ws.Range("F19").Value = 1.041642347744262
$ws.Range("I19").Value = 1.032803272815795
$ws.Range("J19").Value = 1.030085448854999
$ws.Range("K19").Value = 1.036572407038849
$ws.Range("L19").Value = 1.031046001466803
$ws.Range("M19").Value = 1.04493277612805
$ws.Range("N19").Value = 1.014034658815113

# Row 20
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.023617763044015
$ws.Range("D20").Value = 1.033022695959613
$ws.Range("E20").Value = 1.027415574917021
$ws.Range("F20").Value = 1.0413638161382
$ws.Range("I20").Value = 1.032755058169594
$ws.Range("J20").Value = 1.029898417760068
$ws.Range("K20").Value = 1.036406670676936
$ws.Range("L20").Value = 1.030819381299344
$ws.Range("M20").Value = 1.044718724427068
$ws.Range("N20").Value = 1.013971777724106

# Row 21
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.022606753716227
$ws.Range("D21").Value = 1.032273433882735
$ws.Range("E21").Value = 1.026467204914042
$ws.Range("F21").Value = 1.040459363182623
$ws.Range("I21").Value = 1.03259671487239
$ws.Range("J21").Value = 1.029290467416194
$ws.Range("K21").Value = 1.035867496239051
$ws.Range("L21").Value = 1.030083185367824
$ws.Range("M21").Value = 1.044022984001874
$ws.Range("N21").Value = 1.013767294679736

# Row 22
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.021971842789575
$ws.Range("D22").Value = 1.031802821151851
$ws.Range("E22").Value = 1.02587206183415
$ws.Range("F22").Value = 1.039891535273113
$ws.Range("I22").Value = 1.032495940465163
$ws.Range("J22").Value = 1.028908312110262
$ws.Range("K22").Value = 1.035528234434984
$ws.Range("L22").Value = 1.029620752351875
$ws.Range("M22").Value = 1.043585676657004
$ws.Range("N22").Value = 1.013638691673904

# Row 23
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.022308349475578
$ws.Range("D23").Value = 1.03205225597616
$ws.Range("E23").Value = 1.026187450316684
$ws.Range("F23").Value = 1.040192471581777
$ws.Range("I23").Value = 1.032549478099622
$ws.Range("J23").Value = 1.029110891208421
$ws.Range("K23").Value = 1.035708107899183
$ws.Range("L23").Value = 1.02986585443145
$ws.Range("M23").Value = 1.043817488509632
$ws.Range("N23").Value = 1.013706869859529

# Row 24
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.023634725378539
$ws.Range("D24").Value = 1.033035265415324
$ws.Range("E24").Value = 1.027431493579021
$ws.Range("F24").Value = 1.041378993473537
$ws.Range("I24").Value = 1.032757692127442
$ws.Range("J24").Value = 1.029908611519456
$ws.Range("K24").Value = 1.036415705483537
$ws.Range("L24").Value = 1.030831731122982
$ws.Range("M24").Value = 1.044730390737791
$ws.Range("N24").Value = 1.01397520525895

# Row 25
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.025177135166707
$ws.Range("D25").Value = 1.034178017170867
$ws.Range("E25").Value = 1.028880021375857
$ws.Range("F25").Value = 1.042759470822809
$ws.Range("I25").Value = 1.032993989718563
$ws.Range("J25").Value = 1.030834661608462
$ws.Range("K25").Value = 1.037235652625442
$ws.Range("L25").Value = 1.031954458182301
$ws.Range("M25").Value = 1.045790287901
$ws.Range("N25").Value = 1.014286419822408
